# Add data for 2022-01-06: the running "through December 28" tally becomes
# "through December 29" and one additional carjacking record (one per
# neighborhood row, across a handful of monthly columns) is folded into
# the pivot counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2021-12-29"

# Update the running-month column header text to match.
$ws.Range("B1").Value = "December 2021 (through December 29)"

# Increment / add the pivot-table counts touched by the new record(s).
$ws.Range("B2").Value = 3
$ws.Range("Z2").Value = 1
$ws.Range("AL4").Value = 2
$ws.Range("BJ5").Value = 1
$ws.Range("Z6").Value = 6
$ws.Range("B7").Value = 10
$ws.Range("Z7").Value = 8
$ws.Range("B8").Value = 6
$ws.Range("Z8").Value = 2
$ws.Range("BJ8").Value = 6
$ws.Range("B15").Value = 3
$ws.Range("AX16").Value = 4
$ws.Range("B20").Value = 1
$ws.Range("N23").Value = 4
$ws.Range("AX23").Value = 3
$ws.Range("AX34").Value = 1
$ws.Range("B37").Value = 4
$ws.Range("AL37").Value = 2
$ws.Range("B38").Value = 4
$ws.Range("N38").Value = 3
$ws.Range("BV38").Value = 3
$ws.Range("B39").Value = 4
$ws.Range("Z40").Value = 3
$ws.Range("N41").Value = 2
$ws.Range("BV41").Value = 1
$ws.Range("AL57").Value = 3
$ws.Range("B62").Value = 1
$ws.Range("AL62").Value = 1
$ws.Range("B92").Value = 2
$ws.Range("B95").Value = 1
$ws.Range("B97").Value = 1
